$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (2023-09-02 -> 2023-09-03, i.e. serial 45171 -> 45172) for every data row
# (rows 2 through 135).
$ws.Range("C2:C135").Value = 45172
